# This script applies the weekly crime data refresh described in the diff:
#  - bump the report "Number" 51 -> 52
#  - shift the reporting week 12/18/2023-12/24/2023 -> 12/25/2023-12/31/2023
#  - update a handful of crime-count / percentage cells on rows 15,16,19,21,24,25,26

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text tweaks (rich-text shared strings) - edit in place via Characters
# so the surrounding run formatting is preserved.
# ---------------------------------------------------------------------------

# A8 = "Volume 30   Number  51" -> "...52"
$ws.Range("A8").Characters(21, 2).Text = "52"

# C9 = "Report Covering the Week  12/18/2023  Through  12/24/2023"
$ws.Range("C9").Characters(27, 10).Text = "12/25/2023"
$ws.Range("C9").Characters(48, 10).Text = "12/31/2023"

# ---------------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------------
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("F15").Value = 1
$ws.Range("F15").NumberFormat = "#,##0"
$ws.Range("I15").Value = 3
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = -50
$ws.Range("M15").Value = -50
$ws.Range("N15").Value = -72.727272727272

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("L16").Value = 35.294117647058
$ws.Range("N16").Value = -88.725490196078

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
# F19 goes from numeric 2 back to the text placeholder "0" (same as D19/C19)
$ws.Range("D19").Copy($ws.Range("F19"))
$ws.Range("H19").Value = -100

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("C21").Value = 1
$ws.Range("C21").NumberFormat = "#,##0"
$ws.Range("F21").Value = 5
$ws.Range("H21").Value = 25
$ws.Range("I21").Value = 88
$ws.Range("K21").Value = 15.789473684210
$ws.Range("L21").Value = 41.935483870967
$ws.Range("M21").Value = -12.871287128712
$ws.Range("N21").Value = -81.473684210526

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
# D24 goes from numeric 2 to the text placeholder "0"
$ws.Range("D19").Copy($ws.Range("D24"))
# E24 goes from numeric -50 to the text placeholder "***.*"
$ws.Range("E19").Copy($ws.Range("E24"))
$ws.Range("F24").Value = 4
$ws.Range("G24").Value = 4
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 43
$ws.Range("K24").Value = 26.470588235294
$ws.Range("L24").Value = 26.470588235294
$ws.Range("M24").Value = -58.653846153846

# ---------------------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("D25").Value = 1
$ws.Range("J25").Value = 34
$ws.Range("K25").Value = 38.235294117647

# ---------------------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 1
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("F26").Value = 1
$ws.Range("F26").NumberFormat = "#,##0"
$ws.Range("I26").Value = 3
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = -62.5
